$d = $word.ActiveDocument

# 1. Add a new paragraph after the "Week 6" paragraph with the same formatting
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newPara.Range.Text = "Every file of week 6 is important"

# 2. Highlight the "Week 6" text in yellow
$find = $d.Content.Find
$find.Text = "Week 6"
$found = $find.Execute()
if ($found) {
    $find.Parent.Font.HighlightColorIndex = 7  # wdYellow
}
